$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.736.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.531.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.58%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.529.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.989.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.526.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.517.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '358.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.26%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0991'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '553.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.23%  '
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  +1.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.356'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.87%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '149.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.561'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0277'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0759'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.44%  '
